# Update cryptos list values (price + volume%) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'25.895.05"
$ws.Range('E2').Value = "'  +0.53%  "
$ws.Range('D3').Value = "'1.640.07"
$ws.Range('E3').Value = "'  +0.98%  "
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = "'  +0.13%  "
$ws.Range('D5').Value = "'215.81"
$ws.Range('E5').Value = "'  +0.53%  "
$ws.Range('D6').Value = "'0.5083"
$ws.Range('E6').Value = "'  +0.21%  "
$ws.Range('D7').Value = "'1.004"
$ws.Range('E7').Value = "'  +0.20%  "
$ws.Range('D8').Value = "'0.2603"
$ws.Range('E8').Value = "'  +1.87%  "
$ws.Range('D9').Value = "'0.06473"
$ws.Range('E9').Value = "'  +1.63%  "
$ws.Range('D10').Value = "'20.26"
$ws.Range('E10').Value = "'  +5.05%  "
$ws.Range('D11').Value = "'0.07822"
$ws.Range('E11').Value = "'  +0.77%  "
$ws.Range('D12').Value = "'1.665.14"
$ws.Range('E12').Value = "'  +2.52%  "
$ws.Range('D13').Value = "'4.267"
$ws.Range('D14').Value = "'1.866.56"
$ws.Range('E14').Value = "'  +1.00%  "
$ws.Range('D15').Value = "'0.5670"
$ws.Range('E15').Value = "'  +2.56%  "
$ws.Range('D16').Value = "'0.0₅7706"
$ws.Range('E16').Value = "'  +2.71%  "
$ws.Range('D17').Value = "'63.54"
$ws.Range('E17').Value = "'  -0.01%  "
$ws.Range('D18').Value = "'25.910.60"
$ws.Range('E18').Value = "'  +0.51%  "
$ws.Range('D19').Value = "'1.003"
$ws.Range('E19').Value = "'  +0.28%  "
$ws.Range('D20').Value = "'194.84"
$ws.Range('E20').Value = "'  +0.62%  "
$ws.Range('D21').Value = "'4.400"
$ws.Range('E21').Value = "'  +0.22%  "
$ws.Range('E22').Value = "'  +2.47%  "
$ws.Range('D23').Value = "'6.217"
$ws.Range('E23').Value = "'  +4.21%  "
$ws.Range('D24').Value = "'1.004"
$ws.Range('E24').Value = "'  +0.21%  "
$ws.Range('D25').Value = "'1.768"
$ws.Range('E25').Value = "'  -5.19%  "
$ws.Range('D26').Value = "'138.09"
$ws.Range('E26').Value = "'  -1.85%  "
$ws.Range('E27').Value = "'  +0.02%  "
$ws.Range('D28').Value = "'6.873"
$ws.Range('D29').Value = "'15.63"
$ws.Range('E29').Value = "'  +1.40%  "
$ws.Range('E30').Value = "'  +1.02%  "
$ws.Range('D31').Value = "'0.05013"
$ws.Range('E31').Value = "'  +3.20%  "
$ws.Range('D32').Value = "'3.322"
$ws.Range('E32').Value = "'  +0.48%  "
$ws.Range('D33').Value = "'3.260"
$ws.Range('E34').Value = "'  +2.21%  "
$ws.Range('D35').Value = "'2.387"
$ws.Range('E35').Value = "'  +1.09%  "
$ws.Range('D36').Value = "'0.9089"
$ws.Range('E36').Value = "'  +1.97%  "
$ws.Range('D37').Value = "'2.581"
$ws.Range('E37').Value = "'  +1.79%  "
$ws.Range('D38').Value = "'0.5539"
$ws.Range('E38').Value = "'  +0.95%  "
$ws.Range('D39').Value = "'1.131.76"
$ws.Range('E39').Value = "'  +0.16%  "
$ws.Range('D40').Value = "'0.01580"
$ws.Range('E40').Value = "'  +1.62%  "
$ws.Range('D41').Value = "'0.9955"
$ws.Range('E41').Value = "'  -0.53%  "
$ws.Range('D42').Value = "'99.87"
$ws.Range('E42').Value = "'  +2.84%  "
$ws.Range('D43').Value = "'5.496"
$ws.Range('E43').Value = "'  -1.20%  "
$ws.Range('D44').Value = "'0.8039"
$ws.Range('E44').Value = "'  +1.31%  "
$ws.Range('E45').Value = "'  -2.04%  "
$ws.Range('D46').Value = "'55.84"
$ws.Range('E46').Value = "'  +2.25%  "
$ws.Range('D47').Value = "'0.4234"
$ws.Range('E47').Value = "'  -4.15%  "
$ws.Range('D48').Value = "'7.718"
$ws.Range('D49').Value = "'0.05048"
$ws.Range('E49').Value = "'  -1.52%  "
$ws.Range('D50').Value = "'1.002"
$ws.Range('E50').Value = "'  +0.29%  "
$ws.Range('D51').Value = "'1.002"
$ws.Range('E51').Value = "'  +0.15%  "
